$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had two rows (25 and 27) that computed a "speed" (1/x) view of the
# "time" rows directly above them (24 and 26 respectively). Those two
# "speed" rows are being removed; the rows below shift up by two, and the
# "Speedup Factor" / "Relative Speedup Factor" rows (originally 28 and 29,
# now 26 and 27) are re-pointed from the deleted rows to the rows that
# remain (24/25 and 26 respectively).

# Delete original row 25 (" 1 / B24" style formulas).
$ws.Rows("25:25").Delete()
# After that delete, the old row 27 (" 1 / B26" style formulas) is now row 26.
$ws.Rows("26:26").Delete()

# Row 26 ("Speedup Factor Per Data Point When Batched For Neural Networks"):
# was " (B27 / B25)" -> now "B24 / B25" (both its former inputs shifted up).
$ws.Range("B26").Formula = "= B24 / B25"
$ws.Range("C26:G26").FormulaR1C1 = "=R[-2]C/R[-1]C"

# Row 27 ("Relative Speedup Factor..."): was referencing row 28 ($B$28) ->
# now references row 26 ($B$26).
$ws.Range("B27").Formula = "= ((B26 - `$B`$26) / `$B`$26) * 100"
$ws.Range("C27:G27").FormulaR1C1 = "=((R[-1]C-R26C2)/R26C2)*100"

# Fix up the hyperlinks: their anchor cells shifted from row 35 up to row 33
# (the row-delete above doesn't auto-move hyperlink anchors), so drop the
# stale ones and re-add them pointing at the new cells.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D33"), "https://devforum.roblox.com/t/neural-network-library-20/869557/126?u=myoriginsworkshop")
$ws.Hyperlinks.Add($ws.Range("E33"), "https://devforum.roblox.com/t/openml-machine-learning/3008664/4?u=myoriginsworkshop")
$ws.Hyperlinks.Add($ws.Range("F33"), "https://devforum.roblox.com/t/xentorch-neural-network-constructor/1201111/24?u=myoriginsworkshop")
$ws.Hyperlinks.Add($ws.Range("G33"), "https://devforum.roblox.com/t/easyml-an-easy-way-to-use-machine-learning-in-your-roblox-games/3110013?u=myoriginsworkshop")
$ws.Hyperlinks.Add($ws.Range("B33"), "https://devforum.roblox.com/t/datapredict-release-121-general-purpose-machine-learning-and-deep-learning-library-learning-ais-generative-ais-and-more/2196446/2?u=myoriginsworkshop")

# Match the new selection recorded in the workbook view.
$ws.Range("D27").Select()
